$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1420.6586
$ws.Range("I112").Value = 1339.4
$ws.Range("J112").Value = 1431.9445
$ws.Range("K112").Value = 4018.2
$ws.Range("L112").Value = 4295.833500000001
$ws.Range("M112").Value = -2910.2
$ws.Range("N112").Value = -6511.833500000001

$ws.Range("H135").Value = 1850
$ws.Range("I135").Value = 1007.1429
$ws.Range("K135").Value = 9064.286100000001
$ws.Range("M135").Value = -6529.286100000001

$ws.Range("H138").Value = 4107.73
$ws.Range("J138").Value = 4760.1807
$ws.Range("L138").Value = 14280.5421
$ws.Range("N138").Value = -24560.5421

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5576987.5
$ws.Range("I32").Value = 6560356
$ws.Range("J32").Value = 4566.6665
$ws.Range("K32").Value = 6560356
$ws.Range("L32").Value = 4566.6665
$ws.Range("M32").Value = -6560069
$ws.Range("N32").Value = -5140.6665

$ws.Range("H124").Value = 24885.8
$ws.Range("J124").Value = 24885.8
$ws.Range("L124").Value = 24885.8
$ws.Range("N124").Value = -34705.8

$ws.Range("H125").Value = 33638.125
$ws.Range("J125").Value = 33638.125
$ws.Range("L125").Value = 33638.125
$ws.Range("N125").Value = -43478.125

$ws.Range("H132").Value = 1264.8572
$ws.Range("I132").Value = 810.5
$ws.Range("K132").Value = 2431.5
$ws.Range("M132").Value = 98.5

$ws.Range("H135").Value = 29571.6
$ws.Range("J135").Value = 29571.6
$ws.Range("L135").Value = 29571.6
$ws.Range("N135").Value = -39711.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 29000
$ws.Range("J63").Value = 29000
$ws.Range("L63").Value = 29000
$ws.Range("N63").Value = -30372

$ws.Range("H66").Value = 29000
$ws.Range("J66").Value = 29000
$ws.Range("L66").Value = 87000
$ws.Range("N66").Value = -93864

$ws.Range("H130").Value = 54980
$ws.Range("J130").Value = 54980
$ws.Range("L130").Value = 54980
$ws.Range("N130").Value = -65020

$ws.Range("H135").Value = 70780
$ws.Range("J135").Value = 70780
$ws.Range("L135").Value = 70780
$ws.Range("N135").Value = -80920

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3542.5217
$ws.Range("I31").Value = 3544.4546
$ws.Range("J31").Value = 3500
$ws.Range("K31").Value = 3544.4546
$ws.Range("L31").Value = 3500
$ws.Range("M31").Value = -3249.4546
$ws.Range("N31").Value = -4090

$ws.Range("H34").Value = 3542.5217
$ws.Range("I34").Value = 3544.4546
$ws.Range("J34").Value = 3500
$ws.Range("K34").Value = 3544.4546
$ws.Range("L34").Value = 3500
$ws.Range("M34").Value = -3342.4546
$ws.Range("N34").Value = -3904

$ws.Range("H58").Value = 7716.533
$ws.Range("I58").Value = 758
$ws.Range("J58").Value = 12355.556
$ws.Range("K58").Value = 758
$ws.Range("L58").Value = 12355.556
$ws.Range("M58").Value = -555
$ws.Range("N58").Value = -12761.556

$ws.Range("H86").Value = 16007.066
$ws.Range("I86").Value = 35949.832
$ws.Range("J86").Value = 2711.889
$ws.Range("K86").Value = 35949.832
$ws.Range("L86").Value = 2711.889
$ws.Range("M86").Value = -34826.832
$ws.Range("N86").Value = -4957.889

$ws.Range("H89").Value = 16007.066
$ws.Range("I89").Value = 35949.832
$ws.Range("J89").Value = 2711.889
$ws.Range("K89").Value = 179749.16
$ws.Range("L89").Value = 13559.445
$ws.Range("M89").Value = -174133.16
$ws.Range("N89").Value = -24791.445

$ws.Range("H132").Value = 2754.8
$ws.Range("I132").Value = 2664.8
$ws.Range("K132").Value = 7994.400000000001
$ws.Range("M132").Value = -5464.400000000001

$ws.Range("H134").Value = 3452.3635
$ws.Range("I134").Value = 2432.5
$ws.Range("K134").Value = 7297.5
$ws.Range("M134").Value = -4762.5

$ws.Range("H136").Value = 7716.533
$ws.Range("I136").Value = 758
$ws.Range("J136").Value = 12355.556
$ws.Range("K136").Value = 2274
$ws.Range("L136").Value = 37066.66800000001
$ws.Range("M136").Value = 276
$ws.Range("N136").Value = -42166.66800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1040
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 800
$ws.Range("K113").Value = 6000
$ws.Range("L113").Value = 2400
$ws.Range("M113").Value = -3830
$ws.Range("N113").Value = -6740

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2396.2222
$ws.Range("I132").Value = 2575.9565
$ws.Range("J132").Value = 2078.2307
$ws.Range("K132").Value = 7727.869499999999
$ws.Range("L132").Value = 6234.6921
$ws.Range("M132").Value = -5197.869499999999
$ws.Range("N132").Value = -11294.6921

$ws.Range("H136").Value = 22750
$ws.Range("J136").Value = 22750
$ws.Range("L136").Value = 68250
$ws.Range("N136").Value = -73350

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 44490
$ws.Range("J127").Value = 44490
$ws.Range("L127").Value = 44490
$ws.Range("N127").Value = -54410

$ws.Range("H132").Value = 2191
$ws.Range("I132").Value = 1924.5555
$ws.Range("K132").Value = 5773.666499999999
$ws.Range("M132").Value = -3243.666499999999

$ws.Range("H136").Value = 4025.353
$ws.Range("I136").Value = 2045.4615
$ws.Range("J136").Value = 10460
$ws.Range("K136").Value = 6136.3845
$ws.Range("L136").Value = 31380
$ws.Range("M136").Value = -3586.3845
$ws.Range("N136").Value = -36480

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2115.4092
$ws.Range("I81").Value = 1819.8
$ws.Range("J81").Value = 2361.75
$ws.Range("K81").Value = 3639.6
$ws.Range("L81").Value = 4723.5
$ws.Range("M81").Value = -2578.6
$ws.Range("N81").Value = -6845.5

$ws.Range("H84").Value = 2115.4092
$ws.Range("I84").Value = 1819.8
$ws.Range("J84").Value = 2361.75
$ws.Range("K84").Value = 18198
$ws.Range("L84").Value = 23617.5
$ws.Range("M84").Value = -12894
$ws.Range("N84").Value = -34225.5

$ws.Range("H132").Value = 2185.92
$ws.Range("I132").Value = 1875.6
$ws.Range("J132").Value = 2651.4
$ws.Range("K132").Value = 5626.799999999999
$ws.Range("L132").Value = 7954.200000000001
$ws.Range("M132").Value = -3096.799999999999
$ws.Range("N132").Value = -13014.2

$ws.Range("H136").Value = 2583.6072
$ws.Range("I136").Value = 2767.1
$ws.Range("J136").Value = 2124.875
$ws.Range("K136").Value = 8301.299999999999
$ws.Range("L136").Value = 6374.625
$ws.Range("M136").Value = -5751.299999999999
$ws.Range("N136").Value = -11474.625
